$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 182, shifting existing rows 182-252 down to 183-253.
$ws.Rows(182).Insert()

# Populate the new row 182 with the new record.
$ws.Cells.Item(182, 1).Value = 3
$ws.Cells.Item(182, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 45215
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 100112010
$ws.Cells.Item(182, 7).Value = "Achicoria"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 40
$ws.Cells.Item(182, 11).Value = 7000
$ws.Cells.Item(182, 12).Value = 7000
$ws.Cells.Item(182, 13).Value = 7000
$ws.Cells.Item(182, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(182, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(182, 16).Value = 438
$ws.Cells.Item(182, 17).Value = 16
$ws.Cells.Item(182, 18).Value = "Hortaliza"
